# Update cryptocurrency price/volume data in cryptos.xlsx
# Commit: Updated cryptos list on Sat Jul  1 16:56:51 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price column values (e.g. '1.000', '48.34') would otherwise be
# auto-parsed as numbers by Excel; force those specific cells to store
# the new value as plain text, matching the source data feed formatting.
$textPriceCells = @('D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D13', 'D14', 'D15', 'D16', 'D18', 'D20', 'D21', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D45', 'D46', 'D47', 'D48', 'D50', 'D51')
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '30.584.78'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '1.919.66'
$ws.Range('E3').Value = '  +1.95%  '
$ws.Range('E4').Value = '  +0.71%  '
$ws.Range('D5').Value = '247.27'
$ws.Range('E5').Value = '  +4.16%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.67%  '
$ws.Range('D7').Value = '0.4718'
$ws.Range('E7').Value = '  +0.99%  '
$ws.Range('D8').Value = '0.2876'
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('D9').Value = '0.06870'
$ws.Range('E9').Value = '  +4.92%  '
$ws.Range('D10').Value = '104.92'
$ws.Range('E10').Value = '  -3.20%  '
$ws.Range('D11').Value = '18.36'
$ws.Range('E11').Value = '  -2.27%  '
$ws.Range('D12').Value = '1.918.10'
$ws.Range('E12').Value = '  +3.51%  '
$ws.Range('D13').Value = '0.07700'
$ws.Range('E13').Value = '  +2.43%  '
$ws.Range('D14').Value = '5.276'
$ws.Range('E14').Value = '  +4.36%  '
$ws.Range('D15').Value = '0.6705'
$ws.Range('E15').Value = '  +5.63%  '
$ws.Range('D16').Value = '289.23'
$ws.Range('E16').Value = '  -10.03%  '
$ws.Range('D17').Value = '30.589.60'
$ws.Range('E17').Value = '  +1.12%  '
$ws.Range('D18').Value = '0.000007612'
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('E19').Value = '  +0.52%  '
$ws.Range('D20').Value = '12.93'
$ws.Range('E20').Value = '  +0.97%  '
$ws.Range('D21').Value = '5.519'
$ws.Range('E21').Value = '  +8.33%  '
$ws.Range('D22').Value = '2.166.62'
$ws.Range('E22').Value = '  +3.65%  '
$ws.Range('E23').Value = '  +0.73%  '
$ws.Range('D24').Value = '6.308'
$ws.Range('E24').Value = '  +0.77%  '
$ws.Range('D25').Value = '9.379'
$ws.Range('E25').Value = '  +1.92%  '
$ws.Range('D26').Value = '168.46'
$ws.Range('E26').Value = '  +2.11%  '
$ws.Range('D27').Value = '21.12'
$ws.Range('E27').Value = '  +4.43%  '
$ws.Range('D28').Value = '2.124'
$ws.Range('E28').Value = '  +6.89%  '
$ws.Range('D29').Value = '0.1068'
$ws.Range('E29').Value = '  -1.76%  '
$ws.Range('D30').Value = '1.393'
$ws.Range('E30').Value = '  +4.34%  '
$ws.Range('D31').Value = '4.177'
$ws.Range('E31').Value = '  +2.90%  '
$ws.Range('D32').Value = '4.090'
$ws.Range('E32').Value = '  +4.76%  '
$ws.Range('D33').Value = '0.05036'
$ws.Range('E33').Value = '  +1.86%  '
$ws.Range('D34').Value = '0.7372'
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('D35').Value = '1.147'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = '0.02070'
$ws.Range('E36').Value = '  +7.25%  '
$ws.Range('E37').Value = '  +1.79%  '
$ws.Range('D38').Value = '2.687'
$ws.Range('E38').Value = '  +0.28%  '
$ws.Range('D39').Value = '2.060'
$ws.Range('E39').Value = '  +3.30%  '
$ws.Range('D40').Value = '111.31'
$ws.Range('E40').Value = '  +3.97%  '
$ws.Range('D41').Value = '0.8793'
$ws.Range('E41').Value = '  +1.56%  '
$ws.Range('D42').Value = '0.4403'
$ws.Range('E42').Value = '  +7.16%  '
$ws.Range('D43').Value = '5.884'
$ws.Range('E43').Value = '  +2.24%  '
$ws.Range('E44').Value = '  +0.72%  '
$ws.Range('D45').Value = '67.14'
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = '7.279'
$ws.Range('E46').Value = '  +1.45%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').Value = '48.34'
$ws.Range('E47').Value = '  +13.95%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '9.264'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('E49').Value = '  +2.35%  '
$ws.Range('D50').Value = '34.88'
$ws.Range('E50').Value = '  +1.78%  '
$ws.Range('D51').Value = '0.4061'
$ws.Range('E51').Value = '  +7.73%  '
